# Generate Report for handback
#
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values for the second
# handback entry (row 3) on both the "zh-cn" and "de-de" status sheets,
# reflecting the freshly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet - row 3 (9c0760a9-f347-488d-8db0-4bdee644a390 file)
$wsZhCn.Range("D3").Value = "2016-01-07 07:47:14"
$wsZhCn.Range("G3").Value = "2016-01-07 07:47:57"

# de-de sheet - row 3 (9c0760a9-f347-488d-8db0-4bdee644a390 file)
$wsDeDe.Range("D3").Value = "2016-01-07 07:47:24"
$wsDeDe.Range("G3").Value = "2016-01-07 07:48:15"

$wb.Save()
